$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("pt_max"), shifting boson..syst3_c right by one.
$ws.Columns("E:E").Insert()

# New header cell
$ws.Range("E1").Value = "pt_max"

# New data values (pt_max = 50) for all data rows
$ws.Range("E2:E7").Value = 50

# Update the active selection to match the saved workbook view
$ws.Range("E12").Select()
